# Refresh the cryptos price/volume table (columns D "Price" and E
# "Volume(1h)") with the latest scraped figures, per the GitHub Actions
# update commit.
#
# Note: several "Price" values look numeric (e.g. "326.57") but are stored
# as plain text in the workbook (some, like "29.546.97" or "1.919.33", have
# multiple dots and can only ever be text). Assigning a bare numeric-looking
# string to a Range.Value lets Excel auto-convert it to a real number, which
# then round-trips through floating point (326.57 -> 326.56999999999999) and
# corrupts the exact text. Prefixing those assignments with a leading
# apostrophe forces Excel to keep them as literal text, matching the source
# formatting exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.546.97"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.926.76"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").Value = "'1.014"
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "'326.57"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").Value = "'0.4813"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'0.08195"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").Value = "1.919.33"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "'6.091"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "'7.301"
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("D15").Value = "'91.58"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("D16").Value = "'0.06898"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "'0.00001038"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "'17.62"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "29.549.37"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "'5.668"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "'11.99"
$ws.Range("E23").Value = "  +2.25%  "
$ws.Range("D24").Value = "'2.182"
$ws.Range("E24").Value = "  -0.36%  "
$ws.Range("D25").Value = "2.148.47"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'155.97"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").Value = "'6.388"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").Value = "'20.03"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "'2.090"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "'120.54"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").Value = "'0.09589"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").Value = "'5.599"
$ws.Range("E33").Value = "  +1.68%  "
$ws.Range("D34").Value = "'3.566"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'1.386"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "'0.06338"
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("D37").Value = "'0.02279"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").Value = "'1.191"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "'0.5946"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "'10.72"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D42").Value = "'7.885"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "'0.1844"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "'2.466"
$ws.Range("E44").Value = "  +4.10%  "
$ws.Range("D45").Value = "'1.244"
$ws.Range("E45").Value = "  -2.66%  "
$ws.Range("D46").Value = "'12.36"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").Value = "'0.07477"
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").Value = "'0.5549"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "'1.974"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").Value = "'117.61"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").Value = "'2.434"
$ws.Range("E51").Value = "  +1.26%  "
